$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, pushing the existing rows 74-152 down to 75-153.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new record.
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 44810
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = 100112043
$ws.Range("G74").Value = "Pepino ensalada"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 100
$ws.Range("K74").Value = 24000
$ws.Range("L74").Value = 25000
$ws.Range("M74").Value = 24500
$ws.Range("N74").Value = "$/caja 60 unidades"
$ws.Range("O74").Value = "Región de Arica y Parinacota"
$ws.Range("P74").Value = 408
$ws.Range("Q74").Value = 60
$ws.Range("R74").Value = "Hortaliza"
